$d = $word.ActiveDocument

# Locate the end of the last existing bullet ("...no offending rows.") so the
# new Update Log entries are appended right after it, before the sectPr.
$anchor = $d.Content
$found = $anchor.Find.Execute("there were no offending rows.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor paragraph not found"
}
$anchor.Collapse(0)

# Nudge the insertion point one character to the left of the absolute end of
# the document's main story so InsertXML appends a new paragraph after the
# existing content instead of splitting/erasing it at the story boundary.
$insertAt = $d.Range($anchor.End - 1, $anchor.End - 1)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:oel="http://schemas.microsoft.com/office/2019/extlst" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh wp14"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">2/13/2024 – Dan added code to refine ages allowed on age data validation </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>tab</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>0-20 for most species now (was up to 40)</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>0-40 for large catfish spp (spp code 114, 115, 401)</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">0-100 for long lived species (spp code </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>302-307, 311, 321, 325 or 402-405</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">5/29/2024 – Dan </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>bug</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> fix for issue found by Madison Mitchell where having blank gear length caused app to crash.  This resulted from section of code that checks gear length and gear effort to be sure they are reasonable </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>values</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Dray had coded this to identify offending rows using a new column (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>NAlength</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> or </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>NAeffort</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">) where “TRUE” indicates missing values that are required, filtered by those rows, then used a column number to delete the new variable once filtered...but now that we are adding </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">Verified.TL and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>Verified.Wr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">, there are 2 </w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">more columns to the left of the column to be deleted, so this was not deleting the right column resulting in error when </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>rbind</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> was used to merge this </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>dataframe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> with one that has gear length or effort values that were present but out of the legal range. </w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:t>Fix:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">I replaced with a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>dplyr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>select(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>) function to drop the column by column name instead of number</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">Also searched code for [, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>or ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>] to try to find/replace all other instances where column numbers instead of names had been used.</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>...found the following and replaced with select() function:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>In code to download validated sample, column numbers were used to rearrange column order</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> as last step before writing </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="000000"/></w:rPr><w:t>file</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:t xml:space="preserve">Code checking for number of individuals (NOI) is not 0 or </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>NA</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:t xml:space="preserve">Block of code creating </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>output$sampleLength</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> that is used to print row numbers when checkbox is checked.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:t xml:space="preserve">Block of code creating </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>output$sample</w:t></w:r><w:r><w:t>Effort</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> that is used to print row numbers when checkbox is checked.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:t>Code evaluating if NOI is 0 or NA in age data (2 locations in code)</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:t>Code evaluating if TL was blank on age data (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NATLAge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) (2 locations in code)</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr><w:r><w:t>Added some code to put up a modal box if an Age data file were uploaded for sampling validation or if a Sampling data file were uploaded for age data validation.  This prevents crashing app if the wrong file type is loaded.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertAt.InsertXML($xml)
Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
